# Auto-generated Excel COM-interop script
# Reproduces the reordering/expansion of the papers table (A1:G45 -> A1:G47)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Preserve existing cell formats we will need later (date format, wrap-text format) ---
# F4 currently carries the workbook style used for "date_read" cells (numFmtId 14 / m/d/yyyy).
$ws.Range("F4").Copy() | Out-Null
$ws.Range("J1").PasteSpecial($xlPasteFormats) | Out-Null
# C15 currently carries the wrap-text style used for one particular author cell.
$ws.Range("C15").Copy() | Out-Null
$ws.Range("J2").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- Clear the existing data rows (keep header row 1 untouched) ---
$ws.Range("A2:G45").Clear() | Out-Null

# --- Write the final table contents (rows 2-47) ---
# Row 2: Scene Memory Transformer for Embodied Agents in Long-Horizon Tasks
$ws.Range("A2").Value = 'Scene Memory Transformer for Embodied Agents in Long-Horizon Tasks'
$ws.Range("B2").Value = 2019
$ws.Range("C2").Value = 'Kuan Fang, Alexander Toshev, Li Fei-Fei, Silvio Savarese'
$ws.Range("D2").Value = 'seqence to seqence'
$ws.Range("E2").Value = 'attention model'
$ws.Range("G2").Value = 'Transformer powered by Google'

# Row 3: Speech recognition using Dynamic Time Warping
$ws.Range("A3").Value = 'Speech recognition using Dynamic Time Warping'
$ws.Range("B3").Value = 2019
$ws.Range("C3").Value = 'Yurika Permanasari, Erwin H. Harahap, Erwin Prayoga Ali'
$ws.Range("D3").Value = 'speech recognition'
$ws.Range("E3").Value = 'dynamic time warping'
$ws.Range("F3").Value = 43908

# Row 4: Time Delay Recurrent Neural Network for Speech Recognition
$ws.Range("A4").Value = 'Time Delay Recurrent Neural Network for Speech Recognition'
$ws.Range("B4").Value = 2019
$ws.Range("C4").Value = 'Boji Liu,  Weibin Zhang,  Xiangming Xu  and  Dongpeng Chen'
$ws.Range("E4").Value = 'recurrent neural network'

# Row 5: Unsupervised Acoustic Segmentation and Clustering using Siamese Networ
$ws.Range("A5").Value = 'Unsupervised Acoustic Segmentation and Clustering using Siamese Network Embeddings'
$ws.Range("B5").Value = 2019
$ws.Range("C5").Value = 'Saurabhchand Bhati, Shekhar Nayak, K. Sri Rama Murty, Najim Dehak'
$ws.Range("D5").Value = 'phoneme segmentation'

# Row 6: PHONEME BASED EMBEDDED SEGMENTAL K-MEANS FOR UNSUPERVISED TERM DISCOVE
$ws.Range("A6").Value = 'PHONEME BASED EMBEDDED SEGMENTAL K-MEANS FOR UNSUPERVISED TERM DISCOVERY'
$ws.Range("B6").Value = 2018
$ws.Range("C6").Value = 'Saurabhchand Bhati , Herman Kamper and K. Sri Rama Murty'
$ws.Range("D6").Value = 'phoneme segmentation'

# Row 7: Unsupervised Word Segmentation from Speech with Attention
$ws.Range("A7").Value = 'Unsupervised Word Segmentation from Speech with Attention'
$ws.Range("B7").Value = 2018
$ws.Range("C7").Value = 'Pierre Godard et al.'
$ws.Range("D7").Value = 'phoneme segmentation'
$ws.Range("E7").Value = 'attention model'

# Row 8: Attention Is All You Need
$ws.Range("A8").Value = 'Attention Is All You Need'
$ws.Range("B8").Value = 2017
$ws.Range("C8").Value = 'Ashish Vaswani et al.'
$ws.Range("D8").Value = 'seqence to seqence'
$ws.Range("E8").Value = 'attention model'
$ws.Range("F8").Value = 43970
$ws.Range("G8").Value = 'Transformer powered by Google'

# Row 9: Convolutional sequence to sequence learning
$ws.Range("A9").Value = 'Convolutional sequence to sequence learning'
$ws.Range("B9").Value = 2017
$ws.Range("C9").Value = 'Jonas Gehring, Michael  Auli, David  Grangier, Denis Yarats , Yann N Dauphin'
$ws.Range("D9").Value = 'classifier'
$ws.Range("E9").Value = 'convolution neural network'

# Row 10: Segment-Based Speech Emotion Recognition Using Recurrent Neural Networ
$ws.Range("A10").Value = 'Segment-Based Speech Emotion Recognition Using Recurrent Neural Networks'
$ws.Range("B10").Value = 2017
$ws.Range("C10").Value = 'Efthymios Tzinis, Alexandros Potamianos'
$ws.Range("E10").Value = 'attention model'

# Row 11: Unsupervised Phoneme Segmentation Based on Main Energy Change for Arab
$ws.Range("A11").Value = 'Unsupervised Phoneme Segmentation Based on Main Energy Change for Arabic Speech'
$ws.Range("B11").Value = 2017
$ws.Range("C11").Value = 'Lachachi, N. '
$ws.Range("D11").Value = 'phoneme segmentation'

# Row 12: Phonemes based Speech Word Segmentation using K-Means
$ws.Range("A12").Value = 'Phonemes based Speech Word Segmentation using K-Means'
$ws.Range("B12").Value = 2016
$ws.Range("C12").Value = 'Abdulhussein M. Abdullah'
$ws.Range("D12").Value = 'phoneme segmentation'
$ws.Range("E12").Value = 'k-means'

# Row 13: Unsupervised Phoneme Segmentation of Previously Unseen Languages
$ws.Range("A13").Value = 'Unsupervised Phoneme Segmentation of Previously Unseen Languages'
$ws.Range("B13").Value = 2016
$ws.Range("C13").Value = 'Marco Vetter'
$ws.Range("D13").Value = 'speech recognition'

# Row 14: A Neural Algorithm of Artistic Style
$ws.Range("A14").Value = 'A Neural Algorithm of Artistic Style'
$ws.Range("B14").Value = 2015
$ws.Range("C14").Value = 'Leon A. Gatys,  Alexander S. Ecker,  Matthias Bethge'
$ws.Range("D14").Value = 'GAN'
$ws.Range("E14").Value = 'convolution neural network'
$ws.Range("G14").Value = 'image style change'

# Row 15: ADAM: A METHOD FOR STOCHASTIC OPTIMIZATION
$ws.Range("A15").Value = 'ADAM: A METHOD FOR STOCHASTIC OPTIMIZATION'
$ws.Range("B15").Value = 2015
$ws.Range("C15").Value = 'Diederik P. Kingma, Jimmy Lei Ba'
$ws.Range("D15").Value = 'optimizer'
$ws.Range("G15").Value = 'Adam optimizer'

# Row 16: Effective Approaches to Attention-based Neural Machine Translation
$ws.Range("A16").Value = 'Effective Approaches to Attention-based Neural Machine Translation'
$ws.Range("B16").Value = 2015
$ws.Range("C16").Value = 'Minh-Thang Luong, Hieu Pham, Christopher D. Manning'
$ws.Range("E16").Value = 'attention model'

# Row 17: Listen, Attend and Spell
$ws.Range("A17").Value = 'Listen, Attend and Spell'
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 'William Chan, Navdeep Jaitly, Quoc V. Le, Oriol Vinyals'
$ws.Range("D17").Value = 'speech recognition'
$ws.Range("E17").Value = 'attention model'
$ws.Range("F17").Value = 43945
$ws.Range("G17").Value = 'powered by Google, attention model with pyramid-encoder and decoder structure'

# Row 18: NEURAL MACHINE TRANSLATION BY JOINTLY LEARNING TO ALIGN AND TRANSLATE
$ws.Range("A18").Value = 'NEURAL MACHINE TRANSLATION BY JOINTLY LEARNING TO ALIGN AND TRANSLATE'
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 'Dzmitry Bahdanau, KyungHyun Cho, Yoshua Bengio'
$ws.Range("E18").Value = 'attention model'
$ws.Range("G18").Value = 'Bahdanau Attention model'

# Row 19: Discriminative Unsupervised Feature Learning with Convolutional Neural
$ws.Range("A19").Value = 'Discriminative Unsupervised Feature Learning with Convolutional Neural Networks'
$ws.Range("B19").Value = 2014
$ws.Range("C19").Value = 'Alexey Dosovitskiy, Jost Tobias Springenberg, Martin Riedmiller and Thomas Brox'
$ws.Range("D19").Value = 'classifier'
$ws.Range("E19").Value = 'convolution neural network, unsupervised'

# Row 20: Discriminative Unsupervised Feature Learning with Convolutional Neural
$ws.Range("A20").Value = 'Discriminative Unsupervised Feature Learning with Convolutional Neural Networks'
$ws.Range("B20").Value = 2014
$ws.Range("C20").Value = 'Alexey Dosovitskiy, Jost Tobias Springenberg, Martin Riedmiller and Thomas Brox'
$ws.Range("E20").Value = 'convolution neural network'

# Row 21: Learning Phrase Representations using RNN Encoder–Decoder
$ws.Range("A21").Value = 'Learning Phrase Representations using RNN Encoder–Decoder'
$ws.Range("B21").Value = 2014
$ws.Range("C21").Value = 'Kyunghyun Cho, Bart van Merrienboer, Caglar Gulcehre, Dzmitry Bahdanau, Fethi Bougares, Holger Schwenk, Yoshua Bengio'
$ws.Range("D21").Value = 'classifier'
$ws.Range("E21").Value = 'recurrent neural network'
$ws.Range("G21").Value = 'GRU'

# Row 22: Sequence to Sequence Learning with Neural Networks
$ws.Range("A22").Value = 'Sequence to Sequence Learning with Neural Networks'
$ws.Range("B22").Value = 2014
$ws.Range("C22").Value = 'Ilya Sutskever, Oriol Vinyals, Quoc V. Le'
$ws.Range("D22").Value = 'classifier'
$ws.Range("E22").Value = 'recurrent neural network'
$ws.Range("G22").Value = 'powered by Google'

# Row 23: Going Deeper with Convolutions
$ws.Range("A23").Value = 'Going Deeper with Convolutions'
$ws.Range("B23").Value = 2014
$ws.Range("C23").Value = 'Christian Szegedy, Wei Liu, Yangqing Jia, Pierre Sermanet, Scott Reed, Dragomir Anguelov, Dumitru Erhan, Vincent Vanhoucke, Andrew Rabinovich'
$ws.Range("D23").Value = 'image recognition'
$ws.Range("E23").Value = 'convolution neural network'
$ws.Range("G23").Value = 'Inception, average pooling'

# Row 24: ImageNet Classification with Deep Convolutional
$ws.Range("A24").Value = 'ImageNet Classification with Deep Convolutional'
$ws.Range("B24").Value = 2012
$ws.Range("C24").Value = 'Alex Krizhevsky, Ilya Sutskever, Geoffrey E. Hinton'
$ws.Range("D24").Value = 'classifier'
$ws.Range("E24").Value = 'convolution neural network'
$ws.Range("G24").Value = 'ImageNet'

# Row 25: Searching and mining trillions of time series subsequences under dynam
$ws.Range("A25").Value = 'Searching and mining trillions of time series subsequences under dynamic time warping'
$ws.Range("B25").Value = 2012
$ws.Range("C25").Value = 'Thanawin  Rakthanmanon, Bilson Jake L Campana, Abdullah  Mueen, Gustavo E A P A Batista, Brandon  Westover, Qiang  Zhu, Jesin  Zakaria, Eamonn John Keogh'
$ws.Range("E25").Value = 'dynamic time warping'

# Row 26: Separation of Voiced and Unvoiced using Zero crossing rate and Energy 
$ws.Range("A26").Value = 'Separation of Voiced and Unvoiced using Zero crossing rate and Energy of the Speech Signal'
$ws.Range("B26").Value = 2008
$ws.Range("C26").Value = 'Bachu R.G., Kopparthi S., Adapa B., Barkana B.D.'
$ws.Range("D26").Value = 'unvoiced judgement'
$ws.Range("E26").Value = 'zero crossing rate & energy'
$ws.Range("F26").Value = 43928
$ws.Range("G26").Value = 'discriminate voiced/unvoiced signal'

# Row 27: Speech Recognition Using Dynamic Time Warping
$ws.Range("A27").Value = 'Speech Recognition Using Dynamic Time Warping'
$ws.Range("B27").Value = 2008
$ws.Range("C27").Value = 'Talal Bin Amin, Iftekhar Mahmood'
$ws.Range("D27").Value = 'speech recognition'
$ws.Range("E27").Value = 'dynamic time warping'

# Row 28: Unsupervised optimal phoneme segmentation-Objectives, algorithm and co
$ws.Range("A28").Value = 'Unsupervised optimal phoneme segmentation-Objectives, algorithm and comparisons'
$ws.Range("B28").Value = 2008
$ws.Range("C28").Value = 'Yu Qiao, Naoya Shimomura, and Nobuaki Minematsu'
$ws.Range("D28").Value = 'optimizer'
$ws.Range("E28").Value = 'unsupervised'
$ws.Range("F28").Value = 43908
$ws.Range("G28").Value = 'ways to evaluate a unsupervised speech segmentation'

# Row 29: On Clustering Multimedia Time Series Data Using K-Means and Dynamic Ti
$ws.Range("A29").Value = 'On Clustering Multimedia Time Series Data Using K-Means and Dynamic Time'
$ws.Range("B29").Value = 2007
$ws.Range("C29").Value = 'Chotirat Ann Ratanamahatana'

# Row 30: Connectionist Temporal Classification: Labelling Unsegmented Sequence 
$ws.Range("A30").Value = 'Connectionist Temporal Classification: Labelling Unsegmented Sequence Data with Recurrent Neural Networks'
$ws.Range("B30").Value = 2006
$ws.Range("C30").Value = 'Alex Graves, Santiago Fernandez, Faustino Gomez, Jürgen Schmidhuber'
$ws.Range("D30").Value = 'speech recognition'
$ws.Range("E30").Value = 'recurrent neural network'
$ws.Range("F30").Value = 43909
$ws.Range("G30").Value = 'segment speech and recognize speech at the same time with the RNN model'

# Row 31: Phoneme_segmentation_of_speech
$ws.Range("A31").Value = 'Phoneme_segmentation_of_speech'
$ws.Range("B31").Value = 2006
$ws.Range("C31").Value = 'Bartosz Ziółko, Suresh Manandhar and Richard C. Wilson'
$ws.Range("D31").Value = 'phoneme segmentation'
$ws.Range("E31").Value = 'discrete wavelet transform'
$ws.Range("F31").Value = 43903

# Row 32: Unsupervised location-based segmentation of multi-party speech
$ws.Range("A32").Value = 'Unsupervised location-based segmentation of multi-party speech'
$ws.Range("B32").Value = 2004
$ws.Range("C32").Value = 'G. Lathoud, I.A. McCowan and J.M. Odobez'
$ws.Range("D32").Value = 'phoneme segmentation'

# Row 33: Automatic segmentation combining an HMM-based approach and spectral bo
$ws.Range("A33").Value = 'Automatic segmentation combining an HMM-based approach and spectral boundary correction'
$ws.Range("B33").Value = 2002
$ws.Range("C33").Value = 'Yeon-Jun Kim, Alistair Conkie'
$ws.Range("D33").Value = 'phoneme segmentation'
$ws.Range("E33").Value = 'hidden markov model, unsupervised'
$ws.Range("F33").Value = 43908

# Row 34: Least Squares Support Vector Machine Classifiers
$ws.Range("A34").Value = 'Least Squares Support Vector Machine Classifiers'
$ws.Range("B34").Value = 1999
$ws.Range("C34").Value = 'J.A.K. SUYKENS and J. VANDEWALLE'
$ws.Range("D34").Value = 'classifier'
$ws.Range("E34").Value = 'support vector machine'

# Row 35: Markovian Models for Sequential Data
$ws.Range("A35").Value = 'Markovian Models for Sequential Data'
$ws.Range("B35").Value = 1999
$ws.Range("C35").Value = 'Yoshua Bengio'
$ws.Range("D35").Value = 'classifier'
$ws.Range("E35").Value = 'hidden markov model'

# Row 36: Gradient-Based Learning Applied to Document Recognition
$ws.Range("A36").Value = 'Gradient-Based Learning Applied to Document Recognition'
$ws.Range("B36").Value = 1998
$ws.Range("C36").Value = 'Y. LeCun ,  L. Bottou,  Y. Bengio ,  P. Haffner'
$ws.Range("D36").Value = 'classifier'
$ws.Range("E36").Value = 'convolution neural network'
$ws.Range("G36").Value = 'provided by LeCun and won the Turing Award'

# Row 37: LONG SHORT-TERM MEMORY
$ws.Range("A37").Value = 'LONG SHORT-TERM MEMORY'
$ws.Range("B37").Value = 1997
$ws.Range("C37").Value = 'Sepp Hochreiter, Jurgen Schmidhuber'
$ws.Range("D37").Value = 'classifier'
$ws.Range("E37").Value = 'recurrent neural network'
$ws.Range("G37").Value = 'LSTM'

# Row 38: Phoneme segmentation of continuous speech using multi-layer perceptron
$ws.Range("A38").Value = 'Phoneme segmentation of continuous speech using multi-layer perceptron'
$ws.Range("B38").Value = 1996
$ws.Range("C38").Value = 'Youngjoo Suh and Youngiik Lee'
$ws.Range("D38").Value = 'phoneme segmentation'
$ws.Range("E38").Value = 'deep neural network'
$ws.Range("F38").Value = 43911

# Row 39: Preliminary Results on Speech Signal Segmentation with Recurrent Neura
$ws.Range("A39").Value = 'Preliminary Results on Speech Signal Segmentation with Recurrent Neural Networks'
$ws.Range("B39").Value = 1995
$ws.Range("C39").Value = 'Antonio J. Rubio and Ronan G. Reilly'
$ws.Range("D39").Value = 'phoneme segmentation'
$ws.Range("E39").Value = 'recurrent neural network'
$ws.Range("F39").Value = 43907

# Row 40: Connectionist probability estimators in HMM speech recognition
$ws.Range("A40").Value = 'Connectionist probability estimators in HMM speech recognition'
$ws.Range("B40").Value = 1994
$ws.Range("C40").Value = 'S. Renals,  N. Morgan,  H. Bourlard,  M. Cohen,  H. Franco'
$ws.Range("D40").Value = 'classifier'
$ws.Range("E40").Value = 'hidden markov model'
$ws.Range("F40").Value = 43924

# Row 41: Using dynamic time warping to find patterns in time series
$ws.Range("A41").Value = 'Using dynamic time warping to find patterns in time series'
$ws.Range("B41").Value = 1994
$ws.Range("C41").Value = 'Donald J. Bemdt and James Clifford'
$ws.Range("D41").Value = 'time series'
$ws.Range("E41").Value = 'dynamic time warping'
$ws.Range("F41").Value = 43908

# Row 42: Backpropagation Applied to Handwritten Zip Code Recognition
$ws.Range("A42").Value = 'Backpropagation Applied to Handwritten Zip Code Recognition'
$ws.Range("B42").Value = 1989
$ws.Range("C42").Value = 'Y. LeCun,  B. Boser,  J. S. Denker,  D. Henderson'
$ws.Range("D42").Value = 'classifier'
$ws.Range("E42").Value = 'convolution neural network'
$ws.Range("F42").Value = 43928
$ws.Range("G42").Value = 'provided by LeCun'

# Row 43: High performance connected digit recognition using hidden Markov model
$ws.Range("A43").Value = 'High performance connected digit recognition using hidden Markov models'
$ws.Range("B43").Value = 1989
$ws.Range("C43").Value = 'L.R. Rabiner, J.G. Wilpon, F.K. Soong'

# Row 44: Phoneme Recognition Using Time-Delay
$ws.Range("A44").Value = 'Phoneme Recognition Using Time-Delay'
$ws.Range("B44").Value = 1989
$ws.Range("C44").Value = 'A. Waibel,  T. Hanazawa,  G. Hinton,  K. Shikano,  K.J. Lang'
$ws.Range("D44").Value = 'phoneme segmentation'
$ws.Range("E44").Value = 'convolution neural network'
$ws.Range("F44").Value = 43929
$ws.Range("G44").Value = 'first convolution neural network with back propagation but no pooling and it''s 1-dimensional'

# Row 45: Neural networks and physical systems with emergent collective
$ws.Range("A45").Value = 'Neural networks and physical systems with emergent collective'
$ws.Range("B45").Value = 1982
$ws.Range("C45").Value = 'John Joseph Hopfield'
$ws.Range("D45").Value = 'classifier'
$ws.Range("E45").Value = 'recurrent neural network'
$ws.Range("G45").Value = 'first recurrent neural network'

# Row 46: Neocognitron: A Self-organizing Neural Network Model for a Mechanism o
$ws.Range("A46").Value = 'Neocognitron: A Self-organizing Neural Network Model for a Mechanism of Pattern Recognition Unaffected by Shift in Position'
$ws.Range("B46").Value = 1980
$ws.Range("C46").Value = 'Kunihiko Fukushima'
$ws.Range("D46").Value = 'classifier'
$ws.Range("E46").Value = 'convolution neural network'
$ws.Range("G46").Value = 'first convolution neural network which is still not using back propagation'

# Row 47: Where the phonemes are_Dealing with ambiguity in acoustic-phonetic rec
$ws.Range("A47").Value = 'Where the phonemes are_Dealing with ambiguity in acoustic-phonetic recognition'
$ws.Range("B47").Value = 1975
$ws.Range("C47").Value = 'Richard Schwartz and John Makhoul'
$ws.Range("D47").Value = 'phoneme segmentation'
$ws.Range("F47").Value = 43911
$ws.Range("G47").Value = 'segment the most common speech pieces and see if it can be split into smaller parts(multi-step segmentation)'

# --- Re-apply the date format to the "date_read" column cells ---
$ws.Range("J1").Copy() | Out-Null
$ws.Range("F3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F26").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F28").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F30").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F31").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F33").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F38").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F39").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F40").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F41").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F42").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F44").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F47").PasteSpecial($xlPasteFormats) | Out-Null

# --- Re-apply the wrap-text format to the relevant author cell ---
$ws.Range("J2").Copy() | Out-Null
$ws.Range("C44").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- Remove the temporary helper cells ---
$ws.Range("J1").Clear() | Out-Null
$ws.Range("J2").Clear() | Out-Null

# --- Resize the table (ListObject) to match the new data extent ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G47")) | Out-Null

# --- Fix up the sheet view: clear the frozen scroll position and update the selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("G47").Select() | Out-Null

